$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SetsEditor- Proc")

# Update the renamed / corrected set identifiers (rows 17, 20, 21 in columns E, F, H)
$ws.Range("E17").Value = "EPV_RfTp"
$ws.Range("F17").Value = "EPV_RfTp"
$ws.Range("H17").Value = "ERSOLPRI*,ERSOLPRC*,ERSOLPRR*"

$ws.Range("E20").Value = "Ebattery_Dist"
$ws.Range("F20").Value = "Ebattery_Dist"

$ws.Range("E21").Value = "Ebattery_Utility"
$ws.Range("F21").Value = "Ebattery_Utility"

# Update the selected range shown in the worksheet view
$ws.Activate()
$ws.Range("F20:F21").Select()
